$d = $word.ActiveDocument

# 1) Reposition the "_GoBack" bookmark from the first numbered item onto the
#    title paragraph (paragraph 3) *before* removing the two leading empty
#    paragraphs, so the bookmark naturally lands at the very start of the
#    document once those paragraphs are gone.
$d.Bookmarks.Item("_GoBack").Delete()
$titlePara = $d.Paragraphs.Item(3)
$titleStart = $titlePara.Range.Duplicate
$titleStart.Collapse(1)
$d.Bookmarks.Add("_GoBack", $titleStart)

# 2) Remove the two empty centered paragraphs that used to precede the title.
$d.Paragraphs.Item(1).Range.Delete()
$d.Paragraphs.Item(1).Range.Delete()

# 3) Merge ", " + "в первую очередь" + "," into a single ", в первую очередь,".
$d.Content.Find.Execute(", в первую очередь,", $true, $false, $false, $false, $false, $true, 1, $false, ", в первую очередь,", 2)

# 4) Drop the "или, как вариант, ... найти в интернете" alternative, leaving
#    just " (Сделать самим)."
$d.Content.Find.Execute(" (Сделать самим или, как вариант, найти в интернете).", $true, $false, $false, $false, $false, $true, 1, $false, " (Сделать самим).", 2)
